$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The years 2004-2009 (rows 2-7) are removed; remaining rows (2010-2019,
# formerly rows 8-17) shift up to become rows 2-11.
$ws.Range("A2:A7").EntireRow.Delete()
